# Rename the status "swatch" emoji + the "noir" color label to their
# new values everywhere they appear in the used range.
#
#   🟥 -> 📕
#   ⬛ -> 📘
#   🟧 -> 📙
#   🟩 -> 📗
#   noir -> bleu

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "🟥"   = "📕"
    "⬛"   = "📘"
    "🟧"   = "📙"
    "🟩"   = "📗"
    "noir" = "bleu"
}

$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow = $firstRow + $used.Rows.Count - 1
$lastCol = $firstCol + $used.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($map.ContainsKey($v)) {
            $cell.Value = $map[$v]
        }
    }
}
